$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the first
#    (Heading1) paragraph.
# ---------------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$firstRange = $firstPara.Range
$firstRange.Collapse(0)
$firstRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaTargetRange = $d.Range($metaPara.Range.Start, $metaPara.Range.End)

$metaDescriptionText = ": Try Blood Moon Wilds slot game for free today and discover exciting werewolf characters, lunar calendar feature, and more. Compatible on all devices."

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>' + $metaDescriptionText + '</w:t></w:r></w:p>'

$metaTargetRange.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2. Remove the duplicated bold "Play Blood Moon Wilds Free Today..."
#    paragraph further down the document (right before the italic prompt
#    paragraph).
# ---------------------------------------------------------------------------
$found = $true
while ($found) {
    $found = $false
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $pText = $p.Range.Text.TrimEnd("`r")
        if ($i -gt 1 -and $pText -eq "Play Blood Moon Wilds Free Today - Exciting Werewolf Slot Game") {
            $delRange = $d.Range($p.Range.Start, $p.Range.End)
            $delRange.Delete()
            $found = $true
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 3. Replace the text of the closing italic paragraph with the new DALLE
#    image-prompt copy, preserving its italic formatting. Scope the search
#    to the final paragraph only, since the earlier meta-description
#    paragraph contains an overlapping phrase.
# ---------------------------------------------------------------------------
$oldPromptText = "Try Blood Moon Wilds slot game for free today and discover exciting werewolf characters, lunar calendar feature, and more. Compatible on all devices."
$newPromptText = "Prompt: Create a feature image fitting the game Blood Moon Wilds. DALLE, please create a cartoon-style feature image for Blood Moon Wilds that showcases a happy Maya warrior wearing glasses. The image should incorporate elements of the eerie bayou surrounding New Orleans, such as a full moon shining in the background and werewolves lurking in the shadows. The Maya warrior should be holding a treasure chest filled with gold coins and precious jewels, to represent the potential for big wins in the game. Be creative and use bold, vibrant colors to make the image stand out and capture the attention of online slot players."

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Find.Execute($oldPromptText, $true, $false, $false, $false, $false,
                              $true, 1, $false, $newPromptText, 2)
